$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Add the new "2022-Q1" sheet between "2021-Q4" and "总计".
#    Duplicating "2021-Q4" gives us the correct header row / column
#    layout (基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比,
#    持有市值(亿元), 仓位排名) together with the existing header/index
#    cell styling, positioned right after "2021-Q4" automatically.
# -----------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ4.Copy($null, $wsQ4)
$wsQ1 = $wb.Worksheets.Item("2021-Q4 (2)")
$wsQ1.Name = "2022-Q1"

# "2022-Q1" only has one fund row, so drop the duplicated second data row.
$wsQ1.Rows.Item(3).Delete()

# Fill in the single fund row for 2022-Q1.
$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("B2").Value = "002345"
$wsQ1.Range("C2").Value = "华夏高端制造灵活配置混合"

$wsQ1.Range("D2").NumberFormat = "@"
$wsQ1.Range("D2").Value = "23.47"

$wsQ1.Range("E2").NumberFormat = "@"
$wsQ1.Range("E2").Value = "93.20"

$wsQ1.Range("F2").NumberFormat = "@"
$wsQ1.Range("F2").Value = "3.34"

$wsQ1.Range("G2").NumberFormat = "@"
$wsQ1.Range("G2").Value = "0.7839"

$wsQ1.Range("H2").Value = 10

# -----------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: add a 2022-Q1 row above the
#    existing 2021-Q4 row (newest quarter first).
# -----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push the existing 2021-Q4 totals row down to row 3 (copy A2's style
# along with it so the index column keeps its formatting).
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 1.38

# Write the new 2022-Q1 totals row in row 2.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.78
